$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Heading paragraph: "Yearly Report" -> "Report" (bold, 14pt run)
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "Report"
$headRange = $d.Range($p1.Range.Start, $p1.Range.Start + 6)
$headRange.Font.Bold = 1
$headRange.Font.Size = 14

# ---------------------------------------------------------------------
# 2) Replace the big "letter" paragraph with the new set of paragraphs.
#    Build everything in one pass using a single-character placeholder
#    ("#") for paragraph breaks, then convert it to real paragraph
#    marks with a single Find/Replace pass. Two placeholders in a row
#    become two consecutive paragraph marks, i.e. one blank paragraph.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2)

$body = @'
Letter from the Chairperson/President##Dear Friends and Supporters,##As we step into the new year, I am pleased to present our achievements and the challenges we faced in January 2023. This period has been a testament to the resilience and dedication of our team and partners, who continue to work tirelessly towards our vision of creating a better everyday life for the many people.##This month, we embarked on several ambitious projects, including the EAT-Lancet 2.0 Commission, The Climate Story Fund, Climate Resilience Roadmap, and Economic Inclusion as a Pathway to Self-reliance. These initiatives are aligned with our mission to fight poverty and climate change, focusing on sustainable livelihoods and environmental stewardship.##The EAT-Lancet 2.0 Commission has already made significant strides by securing half of its funding[x], setting a strong foundation for promoting sustainable diets globally. Similarly, The Climate Story Fund, with a focus on storytelling to combat climate change, has successfully engaged global audiences and secured substantial funding[x].##Our Climate Resilience Roadmap in East Africa has seen remarkable progress, enhancing local farmers' resilience to climate change[x]. Meanwhile, the project on economic inclusion in Kenya and South Sudan is paving the way for refugees to achieve self-reliance through increased employment opportunities[x].##Despite these successes, we face ongoing challenges such as securing additional funding and managing the geopolitical risks that impact our operations. We are continuously working on strategic partnerships and innovative funding solutions to address these issues[x].##I extend my deepest gratitude to all our supporters and partners. Your unwavering support fuels our commitment to driving positive change. As we move forward, we remain dedicated to our goals, drawing strength from our achievements and learning from the challenges we face.##Together, we are making a difference, one project at a time. Thank you for being part of this journey.##Warm regards,##[Name]#Chairperson/President, IKEA Foundation##---##This letter aims to encapsulate the progress and ongoing efforts of the IKEA Foundation in January 2023, acknowledging both the accomplishments and hurdles while setting a positive tone for the future endeavors.
'@

$p2.Range.Text = $body
$d.Content.Find.Execute("#", $false, $false, $false, $false, $false, $true, 1, $false, "^p", 2)

Write-Host "Paragraph count after split: " $d.Paragraphs.Count

# ---------------------------------------------------------------------
# 3) Fix up the "Letter from the Chairperson/President" paragraph:
#    bold the text run (not the paragraph mark).
# ---------------------------------------------------------------------
$letterPara = $d.Paragraphs(2)
$letterStart = $letterPara.Range.Start
$textLen = "Letter from the Chairperson/President".Length

$boldRange = $d.Range($letterStart, $letterStart + $textLen)
$boldRange.Font.Bold = 1

Write-Host "Done"
